$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Write cell values (rows 9-40) ----
$ws.Range("B9").Value = "Right"
$ws.Range("C9").Value = "Wrong"
$ws.Range("D9").Value = "Not Attempt"
$ws.Range("E9").Value = "Max"
$ws.Range("A10").Value = "No."
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 16
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 28
$ws.Range("A11").Value = "Marking"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0
$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = 55
$ws.Range("C12").Value = -16
$ws.Range("E12").Value = "39/140"
$ws.Range("A15").Value = "Student Ans"
$ws.Range("B15").Value = "Correct Ans"
$ws.Range("D15").Value = "Student Ans"
$ws.Range("E15").Value = "Correct Ans"
$ws.Range("A16").Value = "Option A"
$ws.Range("B16").Value = "Option A"
$ws.Range("D16").Value = "Option A"
$ws.Range("E16").Value = "Option A"
$ws.Range("A17").Value = "Option D"
$ws.Range("B17").Value = "Option D"
$ws.Range("D17").Value = "Option A"
$ws.Range("E17").Value = "Option C"
$ws.Range("A18").Value = "Option B"
$ws.Range("B18").Value = "Option B"
$ws.Range("D18").Value = "Option B"
$ws.Range("E18").Value = "Option D"
$ws.Range("A19").Value = "Option C"
$ws.Range("B19").Value = "Option C"
$ws.Range("A20").Value = "Option B"
$ws.Range("B20").Value = "Option B"
$ws.Range("A21").Value = "Option B"
$ws.Range("B21").Value = "Option C"
$ws.Range("A22").Value = "Option A"
$ws.Range("B22").Value = "Option D"
$ws.Range("A23").Value = "Option B"
$ws.Range("B23").Value = "Option D"
$ws.Range("B24").Value = "Option A"
$ws.Range("A25").Value = "Option C"
$ws.Range("B25").Value = "Option A"
$ws.Range("A26").Value = "Option C"
$ws.Range("B26").Value = "Option C"
$ws.Range("A27").Value = "Option B"
$ws.Range("B27").Value = "Option A"
$ws.Range("A28").Value = "Option B"
$ws.Range("B28").Value = "Option D"
$ws.Range("A29").Value = "Option B"
$ws.Range("B29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("B30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("B31").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("B32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("B33").Value = "Option D"
$ws.Range("A34").Value = "Option C"
$ws.Range("B34").Value = "Option B"
$ws.Range("A35").Value = "Option B"
$ws.Range("B35").Value = "Option D"
$ws.Range("A36").Value = "Option D"
$ws.Range("B36").Value = "Option A"
$ws.Range("A37").Value = "Option B"
$ws.Range("B37").Value = "Option A"
$ws.Range("A38").Value = "Option B"
$ws.Range("B38").Value = "Option A"
$ws.Range("A39").Value = "Option C"
$ws.Range("B39").Value = "Option D"
$ws.Range("A40").Value = "Option B"
$ws.Range("B40").Value = "Option D"

# ---- 2. Apply formatting per style class ----
# style class: plain
$rng = $ws.Range("A9:E9")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("A10:A12")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("D10:E11")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("D12")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108

# style class: green
$rng = $ws.Range("B10:B12")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 32768
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("A16:A20")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 32768
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("D16")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 32768
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("A26")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 32768
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("A30:A33")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 32768
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108

# style class: red
$rng = $ws.Range("C10:C12")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 255
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("D17:D18")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 255
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("A21:A25")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 255
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("A27:A29")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 255
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("A34:A40")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 255
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108

# style class: blue
$rng = $ws.Range("E12")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 16711680
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("B16:B40")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 16711680
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("E16:E18")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Color = 16711680
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108

# style class: bold
$rng = $ws.Range("A15:B15")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Bold = $true
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng = $ws.Range("D15:E15")
$rng.Font.Name = "Century"
$rng.Font.Size = 12
$rng.Font.Bold = $true
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108

Write-Output "edit applied"